# 14 April data update (new date column: 2020-04-14, serial 43934)
# Adds one new trailing date column to each of the four sheets and fills in
# that day's reported figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overal Stats": new column AN (was last col AM)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overal Stats")

# Header date cell - copy the date-number-format from the previous header cell
$ws1.Range("AM1").Copy() | Out-Null
$ws1.Range("AN1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("AN1").Value = 43934

$sheet1Values = [ordered]@{
    "AN3"  = 11518
    "AN4"  = 2058
    "AN5"  = 67
    "AN6"  = 518
    "AN8"  = 94
    "AN9"  = 442
    "AN10" = 212
    "AN11" = 230
    "AN15" = 65
    "AN16" = 50
    "AN17" = 15
    "AN18" = 102
    "AN19" = 152
    "AN20" = 403
    "AN23" = 58
    "AN24" = 52
    "AN25" = 6
    "AN26" = 236
    "AN27" = 288
    "AN28" = 395
    "AN31" = 18
    "AN32" = 18
    "AN33" = 0
    "AN34" = 167
    "AN35" = 185
    "AN36" = 36
    "AN39" = 53
    "AN40" = 44
    "AN41" = 9
    "AN42" = 385
    "AN43" = 429
    "AN44" = 9
    "AN47" = 7
    "AN48" = 6
    "AN49" = 0
    "AN50" = 81
    "AN51" = 87
    "AN52" = 8
    "AN54" = 4
    "AN55" = 3
    "AN56" = 1
    "AN57" = 3
    "AN58" = 6
    "AN59" = 0
    "AN60" = 1
    "AN63" = 62
    "AN64" = 258
    "AN65" = 225
    "AN67" = 43
    "AN68" = 23
    "AN69" = 66
    "AN70" = 9
    "AN72" = 28
    "AN73" = 105
    "AN74" = 105
    "AN75" = 2
}

foreach ($addr in $sheet1Values.Keys) {
    $ws1.Range($addr).Value = $sheet1Values[$addr]
}

# These two totals pick up the thousands-separator number format.
$ws1.Range("AN3").NumberFormat = "#,##0"
$ws1.Range("AN4").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# Sheet "Total Cases by Ward": new column O (was last col N)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Total Cases by Ward")

$ws2.Range("N2").Copy() | Out-Null
$ws2.Range("O2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("O2").Value = 43934

$sheet2Values = [ordered]@{
    "O3"  = 230
    "O4"  = 182
    "O5"  = 149
    "O6"  = 326
    "O7"  = 262
    "O8"  = 290
    "O9"  = 296
    "O10" = 237
    "O11" = 86
}

foreach ($addr in $sheet2Values.Keys) {
    $ws2.Range($addr).Value = $sheet2Values[$addr]
}

# ---------------------------------------------------------------------------
# Sheet "Total Cases by Race": new column J (was last col I)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Total Cases by Race")

$ws3.Range("I2").Copy() | Out-Null
$ws3.Range("J2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws3.Range("J2").Value = 43934

$sheet3Values = [ordered]@{
    "J4"  = 2058
    "J5"  = 400
    "J6"  = 388
    "J7"  = 932
    "J8"  = 30
    "J9"  = 6
    "J10" = 1
    "J11" = 275
    "J12" = 26
    "J14" = 501
    "J15" = 299
    "J16" = 1252
    "J17" = 6
}

foreach ($addr in $sheet3Values.Keys) {
    $ws3.Range($addr).Value = $sheet3Values[$addr]
}

# ---------------------------------------------------------------------------
# Sheet "Lives Lost by Race": new column J (was last col I)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Lives Lost by Race")

$ws4.Range("I1").Copy() | Out-Null
$ws4.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws4.Range("J1").Value = 43934

$sheet4Values = [ordered]@{
    "J3" = 67
    "J4" = 2
    "J5" = 51
    "J6" = 7
    "J7" = 7
    "J8" = 0
}

foreach ($addr in $sheet4Values.Keys) {
    $ws4.Range($addr).Value = $sheet4Values[$addr]
}
